$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 height changed (autofit re-render) 45 -> 30
$ws.Rows.Item(5).RowHeight = 30

# New "Tabela Totem" foreign-key row data (row 35 & 36), entered in the same
# order the original author typed them (matches shared-string insertion order)
$ws.Range("B35").Value = "Código do totem"
$ws.Range("C35").Value = "IdTotem"
$ws.Range("F35").Value = "Auto increment e Chave Primária`n"
$ws.Range("F35").WrapText = $true
$ws.Range("G35").Value = "Identificador da tabela totem, cada totem tem seu próprio id."

$ws.Range("B36").Value = "Chave estrangeira referencia empresa"
$ws.Range("C36").Value = "fk_Empresa"

$ws.Range("D35").Value = "Chave primária numérico inteiro"

$ws.Range("F36").Value = "Preenchimento obrigatório"
$ws.Range("G36").Value = "Referêrencia a tabela empresa."

$ws.Range("E35").Value = 4
$ws.Range("D36").Value = "Numérico inteiro"
$ws.Range("E36").Value = 4

$ws.Rows.Item(35).RowHeight = 60
$ws.Rows.Item(36).RowHeight = 30

# Selection / scroll position (best effort)
$ws.Range("J35").Select()
